# Insert a new data row at row 6 (pushes existing rows 6..79 down to 7..80),
# then populate it with the new weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = 11
$ws.Range("B6").Value = "Vega Monumental Concepción"
$ws.Range("C6").Value = "Bíobío"
$ws.Range("D6").Value = 44496
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 100112043
$ws.Range("G6").Value = "Pepino ensalada"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 350
$ws.Range("K6").Value = 6500
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 6786
$ws.Range("N6").Value = "$/caja 60 unidades"
$ws.Range("O6").Value = "Región de Coquimbo"
$ws.Range("P6").Value = 113
$ws.Range("Q6").Value = 60
$ws.Range("R6").Value = "Hortaliza"
